$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12  = -21.40389999999999
    27  = -22.0156
    32  = -21.27640000000001
    36  = -20.32960000000001
    38  = -19.90879999999999
    46  = -21.8182
    54  = -21.8061
    55  = -22.29179999999999
    56  = -22.11460000000001
    67  = -21.38849999999997
    69  = -21.50819999999998
    72  = -21.65619999999999
    83  = -21.62729999999998
    86  = -21.74309999999999
    91  = -20.32629999999998
    93  = -21.48300000000001
    99  = -21.8495
    104 = -21.24289999999999
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
